$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "Gone are the single-channel days ... viewers. But now these
# dynamics ..." -> "Single-channel days ... viewers are gone, now these
# dynamics ..." split across 4 runs with identical (unchanged) rPr.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Gone are the single-channel days of no choice and no integration and consultation to viewers. But now these dynamics are changing, there are millions of watching options to choose from and they are also available to be streamed across various devices and are really getting user friendly.")
$base1 = $r1.Start
$r1.Text = "Single-channel days of no choice and no integration and consultation to viewers are gone, now these dynamics are changing, there are millions of watching options to choose from and they are also available to be streamed across various devices and are really getting user friendly."

# Boundaries (relative to $base1) of the four resulting runs:
#   [0,1)     "S"
#   [1,79)    "ingle-channel days of no choice and no integration and consultation to viewers"
#   [79,90)   " are gone, "
#   [90,280)  "now these dynamics are changing, ... friendly."
$b1a = $base1 + 1
$b1b = $base1 + 79
$b1c = $base1 + 90
$b1d = $base1 + 280

# Toggling Bold on/off forces the run to split at the range boundaries
# without leaving any residual formatting difference behind.
$x = $d.Range($b1a, $b1b)
$x.Bold = $true
$x.Bold = $false

$x = $d.Range($b1b, $b1c)
$x.Bold = $true
$x.Bold = $false

$x = $d.Range($b1c, $b1d)
$x.Bold = $true
$x.Bold = $false

# ---------------------------------------------------------------------
# Edit 2: "... the insight of the data. This project ..." -> "... the
# insight of the dataset. This project ..." ("set" inserted as its own
# run).
#
# NOTE: inserting text anywhere in the run re-normalises (coalesces)
# every same-formatted run within that paragraph, which would also wipe
# out the existing (untouched) run boundaries between "...Power BI",
# ", Excel" and " to get a visual understanding of the data.". So the
# text insertion is done first with a minimal Find/Replace anchored
# tightly on "data. This", and the run boundaries (both the new ones
# the diff wants, and the pre-existing ones that must be restored) are
# (re)materialised afterwards using the Bold on/off toggle, which is a
# pure-formatting operation and does not trigger that coalescing.
# ---------------------------------------------------------------------
$ins = $d.Content
$ins.Find.Execute("data. This project aims apply", $true, $false, $false, $false, $false, $true, 1, $false, "dataset. This project aims apply", 1)

$r2 = $d.Content
$r2.Find.Execute("The objective of the project is to perform data visualization techniques to understand the insight of the dataset. This project aims apply various Business Intelligence tools such as Power BI")
$base2 = $r2.Start
$end2 = $r2.End

# Boundaries (relative to $base2) of the three resulting runs:
#   [0,110)   "The objective of the project is to perform data visualization techniques to understand the insight of the data"
#   [110,113) "set"
#   [113,191) ". This project aims apply various Business Intelligence tools such as Power BI"
$b2a = $base2 + 110
$b2b = $base2 + 113
# Just after $end2 is the pre-existing ", Excel" run (7 characters) that
# must stay split off from both "...Power BI" (just before it) and
# " to get a visual understanding of the data." (just after it).
$b2c = $end2 + 7

$x = $d.Range($b2a, $b2b)
$x.Bold = $true
$x.Bold = $false

$x = $d.Range($b2b, $end2)
$x.Bold = $true
$x.Bold = $false

$x = $d.Range($end2, $b2c)
$x.Bold = $true
$x.Bold = $false

Write-Output "done"
